$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.8388
$ws.Range("A12").Value = -21.57220000000002
$ws.Range("E14").Value = 16.81490000000001
$ws.Range("E26").Value = 16.33549999999999
$ws.Range("E31").Value = 16.60409999999999
$ws.Range("A32").Value = -21.15869999999999
$ws.Range("E35").Value = 16.3619
$ws.Range("A36").Value = -20.0041
$ws.Range("E37").Value = 16.56640000000002
$ws.Range("A38").Value = -19.41229999999999
$ws.Range("E45").Value = 16.5784
$ws.Range("A46").Value = -21.7013
$ws.Range("A54").Value = -21.96529999999998
$ws.Range("A55").Value = -22.4825
$ws.Range("E57").Value = 16.65500000000001
$ws.Range("A67").Value = -21.42489999999998
$ws.Range("A69").Value = -21.54829999999998
$ws.Range("A72").Value = -21.9451
$ws.Range("A91").Value = -21.42900000000002
$ws.Range("A99").Value = -20.03499999999998
$ws.Range("E100").Value = 16.375
$ws.Range("E102").Value = 16.56679999999999
